# "Update of PS1 calculation modeled values"
# Row 4 (PS1 site) had its measured flux/derived inputs revised; update the
# three underlying measured cells (E4, F4, G4) - dependent formulas
# (H4:K4, and the G15 average) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 0.426
$ws.Range("F4").Value = 5.77
$ws.Range("G4").Value = 78.7

# Leave the selection on G4, matching the saved view state.
$ws.Range("G4").Select()
